$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume cells keep their original plain-text representation
# (these columns store values like "1.001" or "0.000007718" as text, not numbers,
# so force text format before assigning to avoid Excel auto-converting them).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.745.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.948.29"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.60"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4839"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2955"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06833"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "112.67"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.54"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.947.74"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.560"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07666"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6928"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.06"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +8.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.725.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.703"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007718"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.197.34"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.566"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.791"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.07"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.184"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1091"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.437"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.731"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +17.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.429"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05072"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7774"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.163"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02073"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.702"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.045"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.92"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4460"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8753"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.958"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.003"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.393"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.519"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.90"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1253"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.62"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.2552"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.74%  "
